$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(47,48,49,51,53,54,56,58,59,62,63,64,66,71,72,73,74,75,77,78,82,87,88,90,91,92,94,96,99,101,106,112,113,114,115)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "T"
}
